$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 26 de Julio de 2020 a las 00:24"

# Update country names that changed rank/position (column A)
$ws.Range("A9").Value = "Peru"  # was Mexico
$ws.Range("A10").Value = "Mexico"  # was Peru
$ws.Range("A47").Value = "Guatemala"  # was Rumania
$ws.Range("A48").Value = "Rumania"  # was Guatemala
$ws.Range("A80").Value = "Bulgaria"  # was Estado de Palestina
$ws.Range("A81").Value = "Estado de Palestina"  # was Bulgaria
$ws.Range("A109").Value = "Maldivas"  # was Somalia
$ws.Range("A110").Value = "Somalia"  # was Maldivas
$ws.Range("A123").Value = "Suazilandia"  # was Eslovaquia
$ws.Range("A124").Value = "Eslovaquia"  # was Suazilandia
$ws.Range("A162").Value = "Lesoto"  # was Vietnam
$ws.Range("A163").Value = "Vietnam"  # was Lesoto
$ws.Range("A174").Value = "Camboya"  # was Gambia
$ws.Range("A175").Value = "Gambia"  # was Guadalupe
$ws.Range("A176").Value = "Guadalupe"  # was Islas Caimanes
$ws.Range("A177").Value = "Islas Caimanes"  # was Camboya
$ws.Range("A194").Value = "San Martin (Parte Francesa)"  # was Belice
$ws.Range("A195").Value = "Belice"  # was San Martin (Parte Francesa)

# Update numeric statistics cells (columns B-H)
# Row 4
$ws.Range("B4").Value = 4308855
$ws.Range("C4").Value = 60528
$ws.Range("D4").Value = 2053341
$ws.Range("E4").Value = 2106194
$ws.Range("G4").Value = 830
$ws.Range("H4").Value = 149320
# Row 5
$ws.Range("B5").Value = 2394513
$ws.Range("C5").Value = 46313
$ws.Range("E5").Value = 715783
$ws.Range("G5").Value = 1064
$ws.Range("H5").Value = 86449
# Row 9
$ws.Range("B9").Value = 379884
$ws.Range("C9").Value = 3923
$ws.Range("D9").Value = 263130
$ws.Range("E9").Value = 98724
$ws.Range("G9").Value = 187
$ws.Range("H9").Value = 18030
# Row 10
$ws.Range("B10").Value = 378285
$ws.Range("C10").Value = 7573
$ws.Range("D10").Value = 242692
$ws.Range("E10").Value = 92948
$ws.Range("G10").Value = 737
$ws.Range("H10").Value = 42645
# Row 47
$ws.Range("B47").Value = 44492
$ws.Range("C47").Value = 1209
$ws.Range("D47").Value = 31045
$ws.Range("E47").Value = 11748
$ws.Range("G47").Value = 30
$ws.Range("H47").Value = 1699
# Row 48
$ws.Range("B48").Value = 43678
$ws.Range("C48").Value = 1284
$ws.Range("D48").Value = 25373
$ws.Range("E48").Value = 16140
$ws.Range("G48").Value = 15
$ws.Range("H48").Value = 2165
# Row 70
$ws.Range("B70").Value = 15494
$ws.Range("C70").Value = 241
$ws.Range("D70").Value = 9880
$ws.Range("E70").Value = 5520
# Row 80
$ws.Range("B80").Value = 10312
$ws.Range("C80").Value = 189
$ws.Range("D80").Value = 5306
$ws.Range("E80").Value = 4668
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 338
# Row 81
$ws.Range("B81").Value = 10306
$ws.Range("C81").Value = 213
$ws.Range("D81").Value = 3282
$ws.Range("E81").Value = 6949
$ws.Range("G81").Value = 5
$ws.Range("H81").Value = 75
# Row 109
$ws.Range("B109").Value = 3252
$ws.Range("C109").Value = 77
$ws.Range("D109").Value = 2498
$ws.Range("E109").Value = 739
$ws.Range("H109").Value = 15
# Row 110
$ws.Range("B110").Value = 3178
$ws.Range("C110").Value = 7
$ws.Range("D110").Value = 1521
$ws.Range("E110").Value = 1564
$ws.Range("H110").Value = 93
# Row 120
$ws.Range("B120").Value = 2434
$ws.Range("C120").Value = 138
$ws.Range("D120").Value = 518
$ws.Range("E120").Value = 1882
$ws.Range("G120").Value = 2
$ws.Range("H120").Value = 34
# Row 123
$ws.Range("B123").Value = 2142
$ws.Range("C123").Value = 69
$ws.Range("D123").Value = 940
$ws.Range("E123").Value = 1174
# Row 124
$ws.Range("B124").Value = 2141
$ws.Range("C124").Value = 23
$ws.Range("D124").Value = 1577
$ws.Range("E124").Value = 536
# Row 131
$ws.Range("B131").Value = 1752
$ws.Range("C131").Value = 23
$ws.Range("D131").Value = 907
$ws.Range("E131").Value = 840
# Row 149
$ws.Range("D149").Value = 810
$ws.Range("E149").Value = 30
# Row 152
$ws.Range("B152").Value = 853
$ws.Range("C152").Value = 14
$ws.Range("D152").Value = 587
$ws.Range("E152").Value = 249
# Row 162
$ws.Range("B162").Value = 419
$ws.Range("C162").Value = 60
$ws.Range("D162").Value = 69
$ws.Range("E162").Value = 341
$ws.Range("G162").Value = 3
$ws.Range("H162").Value = 9
# Row 163
$ws.Range("B163").Value = 417
$ws.Range("C163").Value = 4
$ws.Range("D163").Value = 365
$ws.Range("E163").Value = 52
$ws.Range("H163").Value = 0
# Row 173
$ws.Range("B173").Value = 263
$ws.Range("C173").Value = 2
$ws.Range("E173").Value = 74
# Row 174
$ws.Range("B174").Value = 225
$ws.Range("C174").Value = 23
$ws.Range("D174").Value = 143
$ws.Range("E174").Value = 82
$ws.Range("H174").Value = 0
# Row 175
$ws.Range("B175").Value = 216
$ws.Range("D175").Value = 60
$ws.Range("E175").Value = 150
$ws.Range("H175").Value = 6
# Row 176
$ws.Range("D176").Value = 176
$ws.Range("E176").Value = 13
$ws.Range("H176").Value = 14
# Row 177
$ws.Range("B177").Value = 203
$ws.Range("D177").Value = 202
$ws.Range("E177").Value = 0
$ws.Range("H177").Value = 1
# Row 184
$ws.Range("D184").Value = 104
$ws.Range("E184").Value = 8
# Row 194
$ws.Range("B194").Value = 49
$ws.Range("C194").Value = 3
$ws.Range("D194").Value = 41
$ws.Range("E194").Value = 5
$ws.Range("H194").Value = 3
# Row 195
$ws.Range("B195").Value = 48
$ws.Range("C195").Value = 1
$ws.Range("D195").Value = 26
$ws.Range("E195").Value = 20
$ws.Range("H195").Value = 2
# Row 217
$ws.Range("B217").Value = 7
$ws.Range("C217").Value = 1
$ws.Range("E217").Value = 1
